$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# --- 1. Rename header labels: _old -> _FV2210, _new -> _FV2304 ---
$oldSuffixCols = @("A","B","C","D","E","F","G","H","I","J")
$newSuffixCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2210")
}
foreach ($col in $newSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2304")
}

# --- 2. Add a Table over the used range ---
$range = $ws.Range("A1:U72")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"

# --- 3. Freeze top row (pane) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
